# Update the "Förändrad" (changed) date column for rows 2-7 on the
# "Avverkningsanmälningar" sheet from 2023-10-25 (45224) to 2023-11-03 (45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45233
}
